$d = $word.ActiveDocument

$ptXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">A microestrutura dos materiais. Sistemas e reticulados cristalinos, grupos espaciais e simetria, tipos mais comuns de estruturas cristalinas. Projeção estereográfica. Direção do feixe difratado e a lei de Bragg. Intensidade do feixe difratado. Métodos de difração de raios X. </w:t><w:br/><w:t>Preparação materialográfica de amostras: corte, embutimento, lixamento e polimento. Técnicas de ataque químico para revelação de fases. Fundamentos de materialografia quantitativa. Microscopia óptica. Técnicas de microscopia eletrônica: varredura e transmissão. Análise química de microrregiões: espectroscopia de energia dispersiva. Técnicas de análise térmica: análise térmica diferencial, calorimetria exploratória diferencial e análise termogravimétrica.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$enXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>The microstructure of materials. Crystal lattices and systems, space groups and symmetry, most common types of crystal structures. Stereographic projection. Direction of the diffracted beam and Bragg's law. Intensity of the diffracted beam. Methods of X-ray diffraction.</w:t><w:br/><w:t>Materialographic sample preparation: cutting, embedding, sanding and polishing. Chemical etching techniques to reveal phases. Fundamentals of quantitative materialography. Optical microscopy. Electron microscopy techniques: scanning and transmission. Chemical analysis of microregions: energy dispersive spectroscopy. Thermal analysis techniques: differential thermal analysis, differential scanning calorimetry and thermogravimetric analysis.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# --- Portuguese paragraph ("Programa" section) ---
$r = $d.Content
$found = $r.Find.Execute("A microestrutura dos materiais.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find Portuguese paragraph" }
$p = $r.Paragraphs(1)
$pr = $p.Range
$pr.InsertXML($ptXml)

# --- English paragraph ("Programa" section, italic) ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("The microstructure of materials.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find English paragraph" }
$p2 = $r2.Paragraphs(1)
$pr2 = $p2.Range
$pr2.InsertXML($enXml)

Write-Host "Done"
